$d = $word.ActiveDocument

# Locate the "LOQ4073: Química Geral II (Requisito fraco)" requirement
# paragraph. The trailing blank paragraph, the "Ver no Jupiter Salvar em
# pdf Salvar em docx" line, and the "© 2020 ..." footer line that follow
# it are leftover scraped page-chrome and should be removed, while the
# final blank paragraph (and the page-break paragraph after it) stay.
$startIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "LOQ4073*") {
        $startIndex = $i + 1
        break
    }
}

if ($startIndex -ne -1) {
    for ($i = 0; $i -lt 3; $i++) {
        $victim = $d.Paragraphs.Item($startIndex)
        $victim.Range.Delete()
    }
}
